$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# The four "low" priority rows (112e67cb, 2d3b79de, 606d9492, eeb061d0) were regenerated
# for handoff: their priority flips from "low" to "ht", and fresh handoff timestamps are
# recorded. Rows 2-3 (0fcde3e5, a33f8897) are untouched.
foreach ($row in 4..7) {
    $wsZh.Range("E$row").Value = "ht"
    $wsDe.Range("E$row").Value = "ht"

    $wsZh.Range("H$row").Value = "2016-09-02 02:37:27"

    $wsOverview.Range("G$row").Value = "2016-09-02 02:37:32"
    $wsDe.Range("H$row").Value = "2016-09-02 02:37:32"
}
